$d = $word.ActiveDocument

$replacements = @(
    @{old = "90% happy. However, the main system integrating supplier was obviously overworked in their projecting contracts. "; new = "Design: 90% happy. However, the main system integrating supplier was obviously overworked in their projecting contracts. "},
    @{old = "The co-operation with Class LR was top professional and very well remembered."; new = "Design: The co-operation with Class LR was top professional and very well remembered."},
    @{old = "Our small purchases for mounting parts was made difficult by payment reputation of ours. Also the delivery lead times of many materials and parts was long."; new = "Design: Our small purchases for mounting parts was made difficult by payment reputation of ours. Also the delivery lead times of many materials and parts was long."},
    @{old = "Generally no complaints."; new = "Design: Generally no complaints."},
    @{old = "Was good."; new = "Design: Was good."},
    @{old = "This system had been succesfully implemented already in NB:s 516 and 517."; new = "Design: This system had been succesfully implemented already in NB:s 516 and 517."}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
